# Apply updated "想去人数" (attendance count) figures to the 广州-漫展信息 workbook,
# matching the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 226
$ws.Range("F5").Value  = 1592
$ws.Range("F7").Value  = 620
$ws.Range("F8").Value  = 132
$ws.Range("F9").Value  = 596
$ws.Range("F10").Value = 51
$ws.Range("F11").Value = 100
$ws.Range("F12").Value = 43
$ws.Range("F13").Value = 161
$ws.Range("F14").Value = 232

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 16
$ws.Range("F12").Value = 196
$ws.Range("F18").Value = 43

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6294
$ws.Range("F3").Value = 781
$ws.Range("F4").Value = 1940
$ws.Range("F5").Value = 25

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 6294
$ws.Range("F3").Value  = 781
$ws.Range("F4").Value  = 1940
$ws.Range("F6").Value  = 25
$ws.Range("F10").Value = 16
$ws.Range("F12").Value = 226
$ws.Range("F16").Value = 1592
$ws.Range("F21").Value = 620
$ws.Range("F22").Value = 196
$ws.Range("F23").Value = 132
$ws.Range("F24").Value = 596
$ws.Range("F25").Value = 51
$ws.Range("F27").Value = 100
$ws.Range("F30").Value = 43
$ws.Range("F31").Value = 161
$ws.Range("F34").Value = 43
$ws.Range("F37").Value = 232
